$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new sample was recorded for 2026/01/18 (same day as the existing row 671),
# so insert a duplicate of row 671 right below it - this carries over the
# matching date/weekday text formatting exactly - then overwrite the two
# columns that actually differ (time-slot hour and ranking).
$ws.Rows(671).Copy()
$ws.Rows(672).Insert()

$ws.Range("C672").Value = 4
$ws.Range("D672").Value = 171

$excel.CutCopyMode = $false
